$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 with the personalized locale share message for Chloe
$ws.Range("B2").Value = "Powerful, creative, and determined. People named Chloe have a deep inner desire for a stable, loving family or community, and a need to work with others and to be appreciated."

# Update selection to B4
$ws.Range("B4").Select()
